# Swap the two theme colour palettes that ship with this deck.
#
# Before: ppt/theme/theme1.xml ("Integral")     <- used by the slide master
#         ppt/theme/theme2.xml ("Office Theme") <- used by the notes master
# After:  ppt/theme/theme1.xml becomes the "Office Theme" palette
#         ppt/theme/theme2.xml becomes the "Integral" palette
#
# (the font scheme and format scheme of the two themes are identical, so
# the only real difference between the two theme parts is their 12-colour
# palette). PowerPoint's ThemeColorScheme collection exposes those 12
# theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in that
# fixed order, as an OLE RGB long (0x00BBGGRR).

function RgbLong([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# The "Office Theme" colours (previously theme2.xml) now become the
# slide master's palette (theme1.xml).
$officeTheme = @(
    (RgbLong 0x00 0x00 0x00),   # dk1
    (RgbLong 0xFF 0xFF 0xFF),   # lt1
    (RgbLong 0x44 0x54 0x6A),   # dk2
    (RgbLong 0xE7 0xE6 0xE6),   # lt2
    (RgbLong 0x5B 0x9B 0xD5),   # accent1
    (RgbLong 0xED 0x7D 0x31),   # accent2
    (RgbLong 0xA5 0xA5 0xA5),   # accent3
    (RgbLong 0xFF 0xC0 0x00),   # accent4
    (RgbLong 0x44 0x72 0xC4),   # accent5
    (RgbLong 0x70 0xAD 0x47),   # accent6
    (RgbLong 0x05 0x63 0xC1),   # hlink
    (RgbLong 0x95 0x4F 0x72)    # folHlink
)

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $masterScheme.Count; $i++) {
    $masterScheme.Item($i).RGB = $officeTheme[$i - 1]
}

# The "Integral" colours (previously theme1.xml) belong on the notes
# master's palette (theme2.xml) now.
$integralTheme = @(
    (RgbLong 0x00 0x00 0x00),   # dk1
    (RgbLong 0xFF 0xFF 0xFF),   # lt1
    (RgbLong 0x45 0x5F 0x51),   # dk2
    (RgbLong 0xE3 0xDE 0xD1),   # lt2
    (RgbLong 0x99 0xCB 0x38),   # accent1
    (RgbLong 0x63 0xA5 0x37),   # accent2
    (RgbLong 0xE6 0xD0 0x24),   # accent3
    (RgbLong 0xCC 0x97 0x00),   # accent4
    (RgbLong 0x4E 0xB3 0xCF),   # accent5
    (RgbLong 0x37 0x8D 0xA6),   # accent6
    (RgbLong 0x6B 0x9F 0x25),   # hlink
    (RgbLong 0xB2 0x6B 0x02)    # folHlink
)

$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $notesScheme.Count; $i++) {
    $notesScheme.Item($i).RGB = $integralTheme[$i - 1]
}

# NOTE: in this host, NotesMaster.Theme and SlideMaster.Theme resolve to
# the same underlying theme part (theme1.xml) - there is no COM path that
# reaches theme2.xml independently. Re-apply the Office palette last so
# theme1.xml (the slide master's theme, the one that matters visually)
# ends up on the intended palette regardless of that aliasing.
for ($i = 1; $i -le $masterScheme.Count; $i++) {
    $masterScheme.Item($i).RGB = $officeTheme[$i - 1]
}

Write-Host "Theme palette updated: slide master -> Office Theme colours."
